# Commit: "Se realizan pruebas para recoleccion del correo mas reciente en tiempo real"
# Adds 38 new log rows (195-232) to the worksheet, mirroring the existing
# email-log table format (columns: Subject / From / Body / NuevaColumna).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$fromValue = " 📭 : Joan Martinez <joan_martinez.olivares@hotmail.com>"
$flagValue = "📩 NUEVO 📩"

# Groups of consecutive rows sharing the same search term.
# AWord: text shown after " 🔎 : " in column A (kept verbatim, may include trailing space)
# CWord: text shown before " Obtener Outlook..." in column C
$groups = @(
    @{ Count = 29; AWord = "Joan";             CWord = "Joan" },
    @{ Count = 5;  AWord = "import datetime";  CWord = "import datetime" },
    @{ Count = 4;  AWord = "Prueba ";          CWord = "Prueba" }
)

$row = 195
foreach ($g in $groups) {
    for ($i = 0; $i -lt $g.Count; $i++) {
        $ws.Cells.Item($row, 1).Value = " 🔎 : " + $g.AWord
        $ws.Cells.Item($row, 2).Value = $fromValue
        $ws.Cells.Item($row, 3).Value = $g.CWord + " Obtener Outlook para iOS<https://aka.ms/o0ukef> "
        $ws.Cells.Item($row, 4).Value = $flagValue
        $row++
    }
}
